$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "legend" table describing the transmission-sign codes used elsewhere
# in the sheet. Appended below the main data table (which currently ends at
# row 61), separated by one blank spacer row (row 62).
# ---------------------------------------------------------------------------

# Row 63 - legend header ("Code" / "Definition"), bold 11pt Times New Roman
# (same font already used e.g. by F2), with a thick rule above and below.
$ws.Range("F2").Copy()
$ws.Range("A63:D63").PasteSpecial(-4122)
$ws.Range("A63:D63").Borders.Item(8).Weight = 4
$ws.Range("A63:D63").Borders.Item(9).Weight = 4

$ws.Range("A63").Value = "Code"
$ws.Range("B63").Value = "Definition"
$ws.Range("B63:D63").HorizontalAlignment = -4108

# Rows 64-68 - legend entries, regular 12pt Times New Roman, left aligned
# only (no explicit vertical centering), the same look used by column A of
# the main table minus the vertical centering.
$ws.Range("A2").Copy()
$ws.Range("A64:D68").PasteSpecial(-4122)
$ws.Range("A64:D68").VerticalAlignment = -4107

# Row 64 gets a thick rule above its B:D block to set the legend body off
# from the header row (mirrors the main table's own row-2 treatment).
$ws.Range("B64:D64").Borders.Item(8).Weight = 4

$ws.Range("A64").Value = "MF"
$ws.Range("B64").Value = "Mesenterial filaments"
$ws.Range("A65").Value = "Mucus"
$ws.Range("B65").Value = "Mucus production"
$ws.Range("A66").Value = "Necrosis"
$ws.Range("B66").Value = "Liquefactive necrosis"
$ws.Range("A67").Value = "Paling"
$ws.Range("B67").Value = "Tissue paling/bleaching"
$ws.Range("A68").Value = "Swelling"
$ws.Range("B68").Value = "Tissue swelling"

# Row 69 - final legend entry (TL / Tissue loss), reuses the existing
# medium-bottom-border style already used to close out the main table
# (e.g. row 61).
$ws.Range("A61").Copy()
$ws.Range("A69:D69").PasteSpecial(-4122)

$ws.Range("A69").Value = "TL"
$ws.Range("B69").Value = "Tissue loss"

# Merge the definition cells (B:D) across every legend row.
$ws.Range("B63:D63").Merge()
$ws.Range("B64:D64").Merge()
$ws.Range("B65:D65").Merge()
$ws.Range("B66:D66").Merge()
$ws.Range("B67:D67").Merge()
$ws.Range("B68:D68").Merge()
$ws.Range("B69:D69").Merge()

Write-Host "legend table added"
